# Applies the FRMediaLMCDAFHIR.xlsx "Mapping de la partie corps" update:
#  1. Bumps the "Date" metadata value.
#  2. Inserts a new top-level row ("FRCDAImageIllustrative" -> "FRMediaDocument",
#     relationship "equivalent") at the top of the "Mapping Table 1" sheet,
#     pushing the existing field-level mapping rows down by one.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "Date" property on the Metadata sheet -------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-28T14:36:08+00:00"

# --- 2. Insert the new top-level mapping row on "Mapping Table 1" ----------
$ws = $wb.Worksheets.Item("Mapping Table 1")

# Row 3 currently holds "FRCDAImageIllustrative.id -> FRMediaDocument.identifier".
# Insert a fresh row above it (shifting rows 3..15 down to 4..16) and fill it
# with the new resource-level mapping.
$ws.Range("A3:E3").EntireRow.Insert()

# Re-apply the same formatting used by the rest of the data rows (the insert
# above leaves the new row with a blank style).
$ws.Range("A4:E4").Copy()
$ws.Range("A3:E3").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A3").Value = "FRCDAImageIllustrative"
$ws.Range("C3").Value = "equivalent"
$ws.Range("D3").Value = "FRMediaDocument"
